$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $text) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
}

# Row 2 - Bitcoin
Set-TextValue "D2" "24.928.22"
Set-TextValue "E2" "  +2.04%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.702.21"
Set-TextValue "E3" "  +0.96%  "

# Row 4 - TetherUSD
Set-TextValue "D4" "1.003"
Set-TextValue "E4" "  +0.15%  "

# Row 5 - BNB
Set-TextValue "D5" "315.69"
Set-TextValue "E5" "  -0.25%  "

# Row 6 - USDC (price unchanged)
Set-TextValue "E6" "  +0.10%  "

# Row 7 - XRP
Set-TextValue "D7" "0.3954"
Set-TextValue "E7" "  +1.71%  "

# Row 8 - Cardano
Set-TextValue "D8" "0.4025"
Set-TextValue "E8" "  +0.24%  "

# Row 9 - Polygon
Set-TextValue "D9" "1.475"
Set-TextValue "E9" "  -0.67%  "

# Row 10 - OKB
Set-TextValue "D10" "52.69"
Set-TextValue "E10" "  +0.39%  "

# Row 11 - BinanceUSD
Set-TextValue "D11" "1.004"
Set-TextValue "E11" "  +0.28%  "

# Row 12 - Dogecoin
Set-TextValue "D12" "0.08809"
Set-TextValue "E12" "  +0.57%  "

# Row 13 - Solana
Set-TextValue "D13" "26.08"
Set-TextValue "E13" "  +0.41%  "

# Row 14 - Polkadot
Set-TextValue "D14" "7.460"
Set-TextValue "E14" "  -0.67%  "

# Row 15 - now ShibaInu (was Chainlink)
Set-TextValue "B15" "ShibaInu"
Set-TextValue "C15" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D15" "0.00001355"
Set-TextValue "E15" "  +0.82%  "

# Row 16 - now Chainlink (was ShibaInu)
Set-TextValue "B16" "Chainlink"
Set-TextValue "C16" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D16" "7.970"
Set-TextValue "E16" "  -0.89%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "1.711.81"
Set-TextValue "E17" "  +2.31%  "

# Row 18 - Litecoin
Set-TextValue "D18" "96.21"
Set-TextValue "E18" "  -1.73%  "

# Row 19 - TRON
Set-TextValue "D19" "0.07175"
Set-TextValue "E19" "  -0.66%  "

# Row 20 - Avalanche
Set-TextValue "D20" "20.55"
Set-TextValue "E20" "  +4.18%  "

# Row 21 - Uniswap
Set-TextValue "D21" "7.342"
Set-TextValue "E21" "  +1.01%  "

# Row 22 - Dai
Set-TextValue "D22" "1.001"
Set-TextValue "E22" "  +0.02%  "

# Row 23 - Cosmos
Set-TextValue "D23" "14.44"
Set-TextValue "E23" "  +1.83%  "

# Row 24 - WrappedBTC
Set-TextValue "D24" "24.927.10"
Set-TextValue "E24" "  +2.02%  "

# Row 25 - now Toncoin (was LidoDAOToken)
Set-TextValue "B25" "Toncoin"
Set-TextValue "C25" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D25" "2.352"
Set-TextValue "E25" "  +0.46%  "

# Row 26 - now LidoDAOToken (was Toncoin)
Set-TextValue "B26" "LidoDAOToken"
Set-TextValue "C26" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D26" "2.964"
Set-TextValue "E26" "  -1.99%  "

# Row 27 - EthereumClassic
Set-TextValue "D27" "23.60"
Set-TextValue "E27" "  +4.42%  "

# Row 28 - HuobiToken
Set-TextValue "D28" "6.177"
Set-TextValue "E28" "  +15.21%  "

# Row 29 - Monero
Set-TextValue "D29" "161.68"
Set-TextValue "E29" "  -3.64%  "

# Row 30 - BitcoinCash
Set-TextValue "D30" "149.91"
Set-TextValue "E30" "  +8.30%  "

# Row 31 - Filecoin
Set-TextValue "D31" "8.345"
Set-TextValue "E31" "  -3.52%  "

# Row 32 - WEMIXTOKEN
Set-TextValue "D32" "2.665"
Set-TextValue "E32" "  +33.93%  "

# Row 33 - WrappedliquidstakedEther2.0
Set-TextValue "D33" "1.900.20"
Set-TextValue "E33" "  +2.39%  "

# Row 34 - Hedera
Set-TextValue "D34" "0.08551"
Set-TextValue "E34" "  -2.26%  "

# Row 35 - VeChain
Set-TextValue "D35" "0.03154"
Set-TextValue "E35" "  +4.42%  "

# Row 36 - ImmutableX
Set-TextValue "D36" "1.052"
Set-TextValue "E36" "  +0.55%  "

# Row 37 - InternetComputer(DFINITY)
Set-TextValue "D37" "7.210"
Set-TextValue "E37" "  -1.87%  "

# Row 38 - Algorand
Set-TextValue "D38" "0.2856"
Set-TextValue "E38" "  +3.26%  "

# Row 39 - now FraxShare (was Stellar)
Set-TextValue "B39" "FraxShare"
Set-TextValue "C39" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D39" "10.88"
Set-TextValue "E39" "  +0.63%  "

# Row 40 - now Stellar (was FraxShare)
Set-TextValue "B40" "Stellar"
Set-TextValue "C40" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D40" "0.09549"
Set-TextValue "E40" "  +4.48%  "

# Row 41 - TheSandbox
Set-TextValue "D41" "0.8259"
Set-TextValue "E41" "  +3.14%  "

# Row 42 - Aptos
Set-TextValue "D42" "13.95"
Set-TextValue "E42" "  -1.00%  "

# Row 43 - TrustWalletToken
Set-TextValue "D43" "1.482"
Set-TextValue "E43" "  +0.67%  "

# Row 44 - EnergySwap
Set-TextValue "D44" "17.50"
Set-TextValue "E44" "  -1.54%  "

# Row 45 - NEARProtocol
Set-TextValue "D45" "2.702"
Set-TextValue "E45" "  +2.97%  "

# Row 46 - Decentraland
Set-TextValue "D46" "0.7379"
Set-TextValue "E46" "  +1.99%  "

# Row 47 - PancakeSwap
Set-TextValue "D47" "4.245"
Set-TextValue "E47" "  -0.84%  "

# Row 48 - Flow
Set-TextValue "D48" "1.415"
Set-TextValue "E48" "  +0.70%  "

# Row 49 - Cronos
Set-TextValue "D49" "0.08762"
Set-TextValue "E49" "  +8.73%  "

# Row 50 - Frax (price unchanged)
Set-TextValue "E50" "  +0.14%  "

# Row 51 - Quant
Set-TextValue "D51" "139.03"
Set-TextValue "E51" "  -0.21%  "
